$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where only D (price) and E (volume) change
$simpleRows = @(
    @{Row=2; DNew='29.495.81'; ENew='  +0.91%  '}
    @{Row=3; DNew='1.881.79'; ENew='  +1.43%  '}
    @{Row=4; DNew='0.9995'; ENew='  -0.16%  '}
    @{Row=5; DNew='0.7163'; ENew='  +2.21%  '}
    @{Row=6; DNew='242.31'; ENew='  +1.92%  '}
    @{Row=7; DNew='0.9997'; ENew='  -0.16%  '}
    @{Row=8; DNew='0.07905'; ENew='  -0.19%  '}
    @{Row=9; DNew='0.3125'; ENew='  +3.56%  '}
    @{Row=10; DNew='25.33'; ENew='  +7.46%  '}
    @{Row=11; DNew='0.08276'; ENew='  +1.09%  '}
    @{Row=12; DNew='0.7323'; ENew='  +3.99%  '}
    @{Row=13; DNew='1.881.94'; ENew='  +1.95%  '}
    @{Row=14; DNew='5.296'; ENew='  +2.17%  '}
    @{Row=16; DNew='29.492.87'; ENew='  +0.97%  '}
    @{Row=17; DNew='5.958'; ENew='  +2.54%  '}
    @{Row=19; DNew='0.000007886'; ENew='  +0.92%  '}
    @{Row=23; DNew='0.9995'; ENew='  -0.17%  '}
    @{Row=24; DNew='0.1618'; ENew='  +14.54%  '}
    @{Row=25; DNew='163.63'; ENew='  +0.68%  '}
    @{Row=26; DNew='9.068'; ENew='  +2.45%  '}
    @{Row=27; DNew='18.39'; ENew='  +1.83%  '}
    @{Row=28; DNew='1.359'; ENew='  -3.23%  '}
    @{Row=29; DNew='1.499'; ENew='  +1.92%  '}
    @{Row=30; DNew='4.396'; ENew='  +1.75%  '}
    @{Row=31; DNew='4.127'; ENew='  +2.91%  '}
    @{Row=32; DNew='0.05281'; ENew='  +2.51%  '}
    @{Row=33; DNew='1.952'; ENew='  +2.12%  '}
    @{Row=35; DNew='0.7293'; ENew='  +2.67%  '}
    @{Row=36; DNew='2.677'; ENew='  -0.15%  '}
    @{Row=37; DNew='0.01875'; ENew='  +1.62%  '}
    @{Row=38; DNew='1.228.36'; ENew='  +6.04%  '}
    @{Row=39; DNew='2.735'; ENew='  +0.91%  '}
    @{Row=40; DNew='0.9142'; ENew='  -1.66%  '}
    @{Row=41; DNew='74.87'; ENew='  +6.75%  '}
    @{Row=42; DNew='6.204'; ENew='  +3.77%  '}
    @{Row=43; DNew='0.9997'; ENew='  -0.16%  '}
    @{Row=44; DNew='102.74'; ENew='  +0.23%  '}
    @{Row=45; DNew='2.043.23'; ENew='  +2.83%  '}
    @{Row=46; DNew='0.5266'; ENew='  -0.53%  '}
    @{Row=49; DNew='9.355'; ENew='  +2.42%  '}
    @{Row=50; DNew='0.4342'; ENew='  +2.27%  '}
    @{Row=51; DNew='7.119'; ENew='  +2.45%  '}
)

foreach ($item in $simpleRows) {
    $dCell = $ws.Cells.Item($item.Row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $item.DNew
    $ws.Cells.Item($item.Row, 5).Value = $item.ENew
}

# Rows where only E (volume) changes
$eOnlyRows = @(
    @{Row=15; ENew='  +2.20%  '}
    @{Row=18; ENew='  +4.95%  '}
    @{Row=20; ENew='  +1.47%  '}
    @{Row=34; ENew='  +3.34%  '}
)

foreach ($item in $eOnlyRows) {
    $ws.Cells.Item($item.Row, 5).Value = $item.ENew
}

# Rows where the coin identity (B/C) swaps along with D/E
$swapRows = @(
    @{Row=21; BNew='Dai'; CNew='https://coinranking.com/coin/MoTuySvg7+dai-dai'; DNew='0.9988'; ENew='  -0.23%  '}
    @{Row=22; BNew='Chainlink'; CNew='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; DNew='7.999'; ENew='  +6.66%  '}
    @{Row=47; BNew='RenderToken'; CNew='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; DNew='1.788'; ENew='  +2.99%  '}
    @{Row=48; BNew='SynthetixNetwork'; CNew='https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'; DNew='2.941'; ENew='  +12.60%  '}
)

foreach ($item in $swapRows) {
    $ws.Cells.Item($item.Row, 2).Value = $item.BNew
    $ws.Cells.Item($item.Row, 3).Value = $item.CNew
    $dCell = $ws.Cells.Item($item.Row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $item.DNew
    $ws.Cells.Item($item.Row, 5).Value = $item.ENew
}
